$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (49) down onto the two
# new rows so the new cells pick up the same cell style ("Bom") used
# throughout the table, without introducing a brand-new style entry.
$ws.Range("A49:E49").Copy()
$ws.Range("A50:E50").PasteSpecial(-4122)
$ws.Range("A49:E49").Copy()
$ws.Range("A51:E51").PasteSpecial(-4122)

# Row 50 - MegaFlare (Greymon damage skill)
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "MegaFlare"
$ws.Range("C50").Value = "DamageSkill(Greymon)"
$ws.Range("D50").Value = 15
$ws.Range("E50").Value = 4

# Row 51 - Dragon'sRoar (Greymon passive skill)
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "Dragon'sRoar"
$ws.Range("C51").Value = "PassiveSkill(Greymon)"
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

$ws.Range("A51:XFD51").Select()
